# Apply the weekly cryptos.xlsx price/volume refresh described by the commit
# 'Updated cryptos list on Sat Nov 18 22:46:48 UTC 2023 with GitHub Actions'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price values (column D) are plain numbers (e.g. "22.13",
# "0.0810"). Excel would normally auto-convert such text into a real number,
# which both changes the cell type and can silently drop significant trailing
# zeros. Mark those specific cells as Text first so the literal string is kept,
# exactly like the other Price cells that already contain inline text.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row-by-row cell updates (Coin / Link / Price / Volume(1h))
$ws.Range("D2").Value = '36.508.61'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '1.957.32'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '243.98'
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("E6").Value = '  -1.03%  '
$ws.Range("D7").Value = '58.71'
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  +2.82%  '
$ws.Range("D10").Value = '0.0810'
$ws.Range("E10").Value = '  -4.63%  '
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("D12").Value = '22.13'
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.828'
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.243.62'
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").Value = '13.69'
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").Value = '1.955.06'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").Value = '36.452.79'
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").Value = '69.67'
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("D20").Value = '0.0₃0856'
$ws.Range("E20").Value = '  -1.31%  '
$ws.Range("D21").Value = '228.33'
$ws.Range("E21").Value = '  -0.82%  '
$ws.Range("D22").Value = '5.05'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").Value = '2.46'
$ws.Range("E24").Value = '  -0.27%  '
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").Value = '  +1.53%  '
$ws.Range("D26").Value = '9.22'
$ws.Range("E26").Value = '  -0.91%  '
$ws.Range("E27").Value = '  -1.58%  '
$ws.Range("D28").Value = '160.34'
$ws.Range("E28").Value = '  -1.11%  '
$ws.Range("D29").Value = '19.42'
$ws.Range("E29").Value = '  -0.66%  '
$ws.Range("E30").Value = '  +0.92%  '
$ws.Range("E31").Value = '  -1.12%  '
$ws.Range("D32").Value = '4.70'
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("D33").Value = '0.0619'
$ws.Range("E33").Value = '  -2.41%  '
$ws.Range("D34").Value = '4.30'
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("E36").Value = '  +2.66%  '
$ws.Range("D37").Value = '3.36'
$ws.Range("E37").Value = '  +9.97%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  -9.47%  '
$ws.Range("D40").Value = '0.0979'
$ws.Range("E40").Value = '  -0.99%  '
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("E42").Value = '  -1.17%  '
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").Value = '16.01'
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("D45").Value = '1.363.56'
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("E46").Value = '  -1.14%  '
$ws.Range("D47").Value = '87.87'
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").Value = '7.12'
$ws.Range("E48").Value = '  -1.48%  '
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").Value = '2.134.98'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '43.65'
$ws.Range("E51").Value = '  -5.00%  '
